$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename players (nicknames -> full real names) and keep their stats,
#     writing rows back out in alphabetical order by the new full name ---
$data = @(
  @("Antonio Rafaschieri", 40, 11, 15, 90, 4, 1),
  @("Antonio Sanzone", 6, 10, 11, 5, 1, 1),
  @("Dave Colaianni", 27, 12, 14, 105, 2, 1),
  @("Domenico Carella", 49, 8, 13, 195, 5, 0),
  @("Fabrizio Fresa", 35, 12, 15, 75, 3, 1),
  @("Giancarlo Tauro", 22, 11, 13, 45, 2, 0),
  @("Giulia Spadafina", 7, 4, 5, 10, 1, 0),
  @("Giuseppe Toto", 23, 9, 11, 80, 2, 0),
  @("Lilly Antonacci", 13, 14, 14, 0, 0, 1),
  @("Luca De Tommasi", 40, 11, 15, 135, 4, 0),
  @("Marco Insabato", 65, 5, 12, 295, 7, 1),
  @("Nicola Cuomo", 52, 6, 11, 210, 5, 0),
  @("Piero Falagario", 7, 6, 7, 20, 1, 0),
  @("Silvio Genchi", 41, 9, 14, 100, 5, 0),
  @("Walter Spadafina", 23, 9, 12, 55, 3, 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}

# --- Widen column A to fit the longer full names (no longer "best fit") ---
$ws.Columns.Item(1).ColumnWidth = 24.6

# --- Update the selected range shown when the sheet is opened ---
[void]$ws.Range("A11").Select()

# --- Stray formatted (underlined) empty cell that shows up near the footer ---
$ws.Range("C20").Font.Underline = 2
